# Extend the "15.1.2" indicator table with two more year columns (2019, 2020).
#
# Row 3 holds the year headers (2013..2018 already present in D3:I3, style
# copied from I3). Row 4 holds the indicator values (6.01 .. 7.38 already
# present in D4:I4, style copied from I4). We add J/K in both rows, copying
# the neighbouring cell's formatting (so borders / fonts / number formats
# stay consistent) before writing the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: year headers -----------------------------------------------
$ws.Range("I3").Copy()
$ws.Range("J3:K3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("J3").Value = 2019
$ws.Range("K3").Value = 2020

# --- Row 4: indicator values ---------------------------------------------
$ws.Range("I4").Copy()
$ws.Range("J4:K4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("J4").Value = 6.18
$ws.Range("K4").Value = 6.18

$excel.CutCopyMode = 0

# --- Selection, matching the author's final cursor position --------------
[void]$ws.Range("G11").Select()
